# modified config of skill
# The "AtkDis" column (H) values are reduced from 3 to 2.5 for every
# data row (rows 2-9). Rows 4-9 previously carried a bordered style
# (matching the look of the header-adjacent rows); after re-entering
# the value the border formatting on those cells is cleared so the
# whole column is visually uniform with rows 2-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the skill data: AtkDis goes from 3 to 2.5 for rows 2-9.
$ws.Range("H2:H9").Value = 2.5

# Clear the border that rows 4-9 had on column H so formatting is
# consistent across the whole column.
$ws.Range("H4:H9").Borders.LineStyle = -4142

# Leave the active selection on the last edited cell.
$ws.Range("H9").Select()
